# Weekly update: a new Espárragos price record (week of 2023-09-05) was
# added to the "Macroferia Regional de Talca" subset, in chronological
# date order, as row 54. Every subsequent existing row shifts down by one.
#
# Strategy: insert a blank row at position 54 (Excel shifts rows 54..97
# down to 55..98 and extends the used range automatically), then populate
# the new row 54 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 54, pushing existing data down.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly record.
$ws.Range("A54").Value = 5
$ws.Range("B54").Value = "Macroferia Regional de Talca"
$ws.Range("C54").Value = "Maule"
$ws.Range("D54").Value = 45174
$ws.Range("E54").Value = 7
$ws.Range("F54").Value = 300000000
$ws.Range("G54").Value = "Espárragos"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 1500
$ws.Range("K54").Value = 2500
$ws.Range("L54").Value = 2500
$ws.Range("M54").Value = 2500
$ws.Range("N54").Value = "$/kilo"
$ws.Range("O54").Value = "Región del Maule"
$ws.Range("P54").Value = 2500
$ws.Range("Q54").Value = 1
$ws.Range("R54").Value = "Hortaliza"
